$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 517
$ws.Cells.Item(5, 6).Value = 90
$ws.Cells.Item(6, 6).Value = 1972
$ws.Cells.Item(7, 6).Value = 6848
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 130
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 890
$ws.Cells.Item(18, 6).Value = 358
$ws.Cells.Item(19, 6).Value = 191
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 2142
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(30, 6).Value = 2910
$ws.Cells.Item(31, 6).Value = 235
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 104
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 975
$ws.Cells.Item(36, 6).Value = 102
$ws.Cells.Item(37, 6).Value = 86
$ws.Cells.Item(40, 6).Value = 259
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 577
$ws.Cells.Item(43, 6).Value = 632
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 35
$ws.Cells.Item(46, 6).Value = 911
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(48, 6).Value = 23
$ws.Cells.Item(50, 6).Value = 147

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 18
$ws.Cells.Item(5, 6).Value = 33
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 62
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(16, 6).Value = 5
$ws.Cells.Item(19, 6).Value = 11
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 108
$ws.Cells.Item(26, 6).Value = 49
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 4
$ws.Cells.Item(30, 6).Value = 1

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 6352

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 517
$ws.Cells.Item(5, 6).Value = 90
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 11560
$ws.Cells.Item(11, 6).Value = 12378
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 1286
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 5313
$ws.Cells.Item(16, 6).Value = 890
$ws.Cells.Item(17, 6).Value = 358
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(22, 6).Value = 320
$ws.Cells.Item(23, 6).Value = 1972
$ws.Cells.Item(24, 6).Value = 985
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(26, 6).Value = 482
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(29, 6).Value = 235
$ws.Cells.Item(30, 6).Value = 1967
$ws.Cells.Item(31, 6).Value = 104
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(37, 6).Value = 975
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 49
$ws.Cells.Item(42, 6).Value = 577
$ws.Cells.Item(44, 6).Value = 136
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 6).Value = 273
$ws.Cells.Item(47, 6).Value = 23
$ws.Cells.Item(48, 6).Value = 4272
$ws.Cells.Item(49, 6).Value = 147
$ws.Cells.Item(50, 6).Value = 0
